$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for 9121c07e-...md row
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-31 07:31:41"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime for 9121c07e-...md row
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-08-31 07:31:30"
$wsZhCn.Range("K3").Value = "2016-08-31 07:32:28"

# de-de sheet: Correspond Handoff Datetime for 9121c07e-...md row
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-08-31 07:32:47"
